$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new columns E:H
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "size"
$ws.Range("G1").Value = "layout"
$ws.Range("H1").Value = "domain_name"

# Data rows 2-19: type, size, layout, domain_name derived from the filename in column A
$data = @(
    @("node-link","xl","hier","drinking"),
    @("node-link","xl","organic","drinking"),
    @("node-link","xl","radial","drinking"),
    @("node-link","med","hier","drinking"),
    @("node-link","med","organic","drinking"),
    @("node-link","med","radial","drinking"),
    @("node-link","xl","hier","sport"),
    @("node-link","xl","organic","sport"),
    @("node-link","xl","radial","sport"),
    @("node-link","med","hier","sport"),
    @("node-link","med","organic","sport"),
    @("node-link","med","radial","sport"),
    @("node-link","xl","hier","student"),
    @("node-link","xl","organic","student"),
    @("node-link","xl","radial","student"),
    @("node-link","med","hier","student"),
    @("node-link","med","organic","student"),
    @("node-link","med","radial","student")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 5).Value = $vals[0]
    $ws.Cells.Item($row, 6).Value = $vals[1]
    $ws.Cells.Item($row, 7).Value = $vals[2]
    $ws.Cells.Item($row, 8).Value = $vals[3]
}

# New column widths for D and E (closest representable values - the host
# quantizes ColumnWidth to 1/6-character steps)
$ws.Columns.Item(4).ColumnWidth = 23.5
$ws.Columns.Item(5).ColumnWidth = 37.666666666666664

# View state: scroll so column B is the leftmost visible, select E5
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E5").Select()
